$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (Price and Volume(1h) columns).
# Values must be kept as literal text (e.g. "311.90", "0.65%") rather than
# being auto-converted to numbers/percentages, so NumberFormat is forced to
# Text ("@") before assigning each value.
$updates = [ordered]@{
    'D2' = '311.90'
    'E2' = '0.65%'
    'D3' = '37.67'
    'E3' = '-0.01%'
    'D4' = '5.137'
    'E4' = '0.62%'
    'D5' = '0.07897'
    'E5' = '0.48%'
    'D6' = '4.409'
    'E6' = '0.90%'
    'D7' = '8.271'
    'E7' = '-0.39%'
    'E8' = '-3.85%'
    'D9' = '0.9257'
    'E9' = '-0.26%'
    'D10' = '0.1215'
    'E10' = '-9.78%'
    'D11' = '0.1933'
    'E11' = '-2.86%'
    'D12' = '0.09134'
    'E12' = '2.17%'
    'D13' = '0.03285'
    'E13' = '-5.23%'
    'D14' = '0.09624'
    'E14' = '-0.76%'
    'D15' = '0.001377'
    'E15' = '-0.81%'
    'D16' = '0.005798'
    'E16' = '-2.82%'
    'D17' = '3.514'
    'E17' = '-2.02%'
    'D18' = '3.097'
    'E18' = '-1.15%'
    'D19' = '0.3398'
    'E19' = '-1.94%'
    'D20' = '5.275'
    'E20' = '5.46%'
    'D21' = '0.1281'
    'E21' = '-1.05%'
    'D22' = '0.2589'
    'E22' = '3.09%'
    'D24' = '0.04359'
    'E24' = '0.50%'
    'D25' = '0.001238'
    'E25' = '0.94%'
    'D26' = '0.004312'
    'E26' = '-5.12%'
    'D27' = '0.0001219'
    'E27' = '-9.70%'
    'D39' = '0.02121'
    'E39' = '-7.74%'
    'D40' = '0.05178'
    'E40' = '2.17%'
    'D41' = '0.007629'
    'E41' = '2.09%'
    'D42' = '0.009140'
    'E42' = '-7.10%'
    'D43' = '0.1362'
    'E43' = '0.26%'
    'D44' = '0.002019'
    'E44' = '1.97%'
    'D45' = '0.008598'
    'E45' = '-2.05%'
    'D46' = '0.00006712'
    'E46' = '-1.69%'
    'E47' = '-0.09%'
    'D48' = '0.001200'
    'E48' = '-7.77%'
    'D49' = '0.002815'
    'E49' = '-6.21%'
    'D50' = '0.00002099'
    'E50' = '-0.09%'
    'D51' = '0.0001999'
    'E51' = '-0.09%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$addr]
}
